$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers: G1 text changes from "FR" to "Avg FR", new header "Avg BW" added at H1
$ws.Range("G1").Value = "Avg FR"
$ws.Range("H1").Value = "Avg BW"

# Full replacement data block for A2:H13
$data = @(
    @(1, 1080, 8.33, 0.32, $null, 67.770475000000005, 109.83050847457601, 468.317288135593),
    @(1, 1080, 16.66, 0.64, $null, 93.222066666666606, 57.745762711864401, 246.227932203389),
    @(1, 720, 8.33, 0.32, $null, 68.371228915662599, 110.186440677966, 469.83498305084697),
    @(1, 720, 16.66, 0.64, $null, 91.800146341463403, 57.711864406779597, 246.083389830508),
    @(2, 1080, 8.33, 0.32, $null, 120.027583333333, 28.932203389830502, 123.366915254237),
    @(2, 1080, 16.66, 0.64, $null, 127.963130434782, 28.677966101694899, 122.282847457627),
    @(2, 720, 8.33, 0.32, $null, 119.45613636363601, 28.864406779661, 123.07783050847399),
    @(2, 720, 16.66, 0.64, $null, 126.154368421052, 28.745762711864401, 122.571932203389),
    @(3, 1080, 8.33, 0.32, $null, 149.988888888888, 19.982758620689602, 85.206482758620595),
    @(3, 1080, 16.66, 0.64, $null, 161.73072727272699, 19.982758620689602, 85.206482758620595),
    @(3, 720, 8.33, 0.32, $null, 152.130764705882, 19.982758620689602, 85.206482758620595),
    @(3, 720, 16.66, 0.64, $null, 155.1001, 19.982758620689602, 85.206482758620595)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# Update selection to match target state
$ws.Range("K9").Select()
